$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 115 (existing rows 115-121 shift down to 117-123).
$ws.Rows("115:116").Insert()

# Common (unchanged) columns shared by every data row in this block.
$mercadoId = 4
$mercado   = "Feria Lagunitas de Puerto Montt"
$region    = "Los Lagos"
$codreg    = 10
$tipo      = "Fruta"
$productoId = 100103
$producto   = "Frutos de hueso (carozo)"
$categoriaId = 100103002
$categoria   = "Ciruela"

# New row 115: Lemon / Primera
$r = 115
$ws.Cells.Item($r,1).Value2 = $mercadoId
$ws.Cells.Item($r,2).Value2 = $mercado
$ws.Cells.Item($r,3).Value2 = $region
$ws.Cells.Item($r,4).Value2 = 44578
$ws.Cells.Item($r,5).Value2 = $codreg
$ws.Cells.Item($r,6).Value2 = $tipo
$ws.Cells.Item($r,7).Value2 = $productoId
$ws.Cells.Item($r,8).Value2 = $producto
$ws.Cells.Item($r,9).Value2 = $categoriaId
$ws.Cells.Item($r,10).Value2 = $categoria
$ws.Cells.Item($r,11).Value2 = "Lemon"
$ws.Cells.Item($r,12).Value2 = "Primera"
$ws.Cells.Item($r,13).Value2 = 200
$ws.Cells.Item($r,14).Value2 = 17000
$ws.Cells.Item($r,15).Value2 = 17500
$ws.Cells.Item($r,16).Value2 = 17250
$ws.Cells.Item($r,17).Value2 = "`$/caja 15 kilos granel"
$ws.Cells.Item($r,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item($r,19).Value2 = 1150
$ws.Cells.Item($r,20).Value2 = 15

# New row 116: Lemon / Segunda
$r = 116
$ws.Cells.Item($r,1).Value2 = $mercadoId
$ws.Cells.Item($r,2).Value2 = $mercado
$ws.Cells.Item($r,3).Value2 = $region
$ws.Cells.Item($r,4).Value2 = 44578
$ws.Cells.Item($r,5).Value2 = $codreg
$ws.Cells.Item($r,6).Value2 = $tipo
$ws.Cells.Item($r,7).Value2 = $productoId
$ws.Cells.Item($r,8).Value2 = $producto
$ws.Cells.Item($r,9).Value2 = $categoriaId
$ws.Cells.Item($r,10).Value2 = $categoria
$ws.Cells.Item($r,11).Value2 = "Lemon"
$ws.Cells.Item($r,12).Value2 = "Segunda"
$ws.Cells.Item($r,13).Value2 = 100
$ws.Cells.Item($r,14).Value2 = 14000
$ws.Cells.Item($r,15).Value2 = 14000
$ws.Cells.Item($r,16).Value2 = 14000
$ws.Cells.Item($r,17).Value2 = "`$/caja 15 kilos granel"
$ws.Cells.Item($r,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item($r,19).Value2 = 933
$ws.Cells.Item($r,20).Value2 = 15
